$wb = $excel.ActiveWorkbook

# --- 1) Metadata sheet: bump "Last Updated" timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 02:07 PM"

# --- 2) Stock List sheet: drop the first two tickers (MIDWESTLTD, CAPTRU-RE1)
#         which shifts every remaining row up by two, then append two new
#         tickers (SMARTWORKS, TRAVELFOOD) at the freed-up bottom rows. ---
$ws = $wb.Worksheets.Item("Stock List")

$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()

$ws.Range("A75").Value = "📋"
$ws.Range("B75").Value = "SMARTWORKS"
$ws.Range("C75").Value = "SMARTWORKS"
$ws.Range("D75").Value = 606.65
$ws.Range("E75").Value = 2.0867
$ws.Range("F75").Value = "N/A"
$ws.Range("G75").Value = "N/A"
$ws.Range("H75").Value = 6931.2448

$ws.Range("A76").Value = "📋"
$ws.Range("B76").Value = "TRAVELFOOD"
$ws.Range("C76").Value = "TRAVELFOOD"
$ws.Range("D76").Value = 1316.3
$ws.Range("E76").Value = 0.1141
$ws.Range("F76").Value = "N/A"
$ws.Range("G76").Value = "N/A"
$ws.Range("H76").Value = 17332.9705
